# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns,
# plus fix the swapped FTXToken/Algorand rows (44/45).
#
# Values must stay plain text cells (inlineStr in the source OOXML), just
# like before the edit. Excel auto-converts a single-dot numeric-looking
# string (e.g. "4.09") to a real number on assignment, so for those cells
# we briefly force the Text number format, assign the value, then clear
# the format again so the cell's style index is left exactly as it was.

function Set-TextCell($addr, $text, $looksNumeric) {
    $cell = $ws.Range($addr)
    if ($looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.ClearFormats()
    } else {
        $cell.Value = $text
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell 'D2' '42.055.32' $false
Set-TextCell 'E2' '  -0.63%  ' $false
Set-TextCell 'D3' '2.232.45' $false
Set-TextCell 'E3' '  -0.60%  ' $false
Set-TextCell 'D5' '250.07' $true
Set-TextCell 'E5' '  +7.06%  ' $false
Set-TextCell 'D6' '0.631' $true
Set-TextCell 'E6' '  -0.41%  ' $false
Set-TextCell 'D7' '72.10' $true
Set-TextCell 'E7' '  +3.43%  ' $false
Set-TextCell 'E8' '  -0.04%  ' $false
Set-TextCell 'D9' '0.596' $true
Set-TextCell 'E9' '  +6.26%  ' $false
Set-TextCell 'D10' '41.37' $true
Set-TextCell 'E10' '  +15.00%  ' $false
Set-TextCell 'D11' '0.0973' $true
Set-TextCell 'E11' '  -1.77%  ' $false
Set-TextCell 'D12' '58.14' $true
Set-TextCell 'E12' '  -0.34%  ' $false
Set-TextCell 'D13' '7.19' $true
Set-TextCell 'E13' '  +6.37%  ' $false
Set-TextCell 'E14' '  -0.61%  ' $false
Set-TextCell 'D15' '2.565.33' $false
Set-TextCell 'E15' '  -0.64%  ' $false
Set-TextCell 'D16' '15.02' $true
Set-TextCell 'E16' '  -0.20%  ' $false
Set-TextCell 'D17' '0.865' $true
Set-TextCell 'E17' '  +0.58%  ' $false
Set-TextCell 'D18' '2.224.87' $false
Set-TextCell 'E18' '  -1.00%  ' $false
Set-TextCell 'D19' '41.933.81' $false
Set-TextCell 'E19' '  -0.71%  ' $false
Set-TextCell 'E20' '  -0.92%  ' $false
Set-TextCell 'D21' '6.24' $true
Set-TextCell 'E21' '  -0.54%  ' $false
Set-TextCell 'D22' '73.04' $true
Set-TextCell 'E22' '  -0.55%  ' $false
Set-TextCell 'D23' '235.87' $true
Set-TextCell 'E23' '  -0.26%  ' $false
Set-TextCell 'E24' '  +8.26%  ' $false
Set-TextCell 'D25' '4.09' $true
Set-TextCell 'E25' '  +11.69%  ' $false
Set-TextCell 'E26' '  -0.06%  ' $false
Set-TextCell 'D27' '2.54' $true
Set-TextCell 'E27' '  +7.09%  ' $false
Set-TextCell 'D28' '10.81' $true
Set-TextCell 'E28' '  +7.86%  ' $false
Set-TextCell 'D29' '171.85' $true
Set-TextCell 'E29' '  +1.54%  ' $false
Set-TextCell 'E30' '  -3.20%  ' $false
Set-TextCell 'D31' '20.79' $true
Set-TextCell 'E31' '  +1.05%  ' $false
Set-TextCell 'E32' '  +3.32%  ' $false
Set-TextCell 'D33' '0.126' $true
Set-TextCell 'E33' '  -0.98%  ' $false
Set-TextCell 'D34' '5.55' $true
Set-TextCell 'E34' '  +3.61%  ' $false
Set-TextCell 'D35' '0.0735' $true
Set-TextCell 'E35' '  +2.04%  ' $false
Set-TextCell 'D36' '4.73' $true
Set-TextCell 'E36' '  +0.12%  ' $false
Set-TextCell 'D37' '26.10' $true
Set-TextCell 'E37' '  +21.22%  ' $false
Set-TextCell 'D38' '3.98' $true
Set-TextCell 'E38' '  +9.77%  ' $false
Set-TextCell 'D39' '0.0301' $true
Set-TextCell 'E39' '  +12.26%  ' $false
Set-TextCell 'E40' '  +1.41%  ' $false
Set-TextCell 'D41' '6.02' $true
Set-TextCell 'E41' '  +1.15%  ' $false
Set-TextCell 'D42' '67.48' $true
Set-TextCell 'E42' '  +1.93%  ' $false
Set-TextCell 'D43' '12.04' $true
Set-TextCell 'E43' '  +20.24%  ' $false
Set-TextCell 'B44' 'FTXToken' $false
Set-TextCell 'C44' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' $false
Set-TextCell 'D44' '5.00' $true
Set-TextCell 'E44' '  +1.98%  ' $false
Set-TextCell 'B45' 'Algorand' $false
Set-TextCell 'C45' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' $false
Set-TextCell 'D45' '0.205' $true
Set-TextCell 'E45' '  +8.65%  ' $false
Set-TextCell 'D46' '8.79' $true
Set-TextCell 'E46' '  -1.85%  ' $false
Set-TextCell 'E47' '  -0.64%  ' $false
Set-TextCell 'D48' '4.67' $true
Set-TextCell 'E48' '  +5.28%  ' $false
Set-TextCell 'E49' '  -0.27%  ' $false
Set-TextCell 'D50' '1.17' $true
Set-TextCell 'E50' '  +7.46%  ' $false
Set-TextCell 'D51' '1.20' $true
Set-TextCell 'E51' '  +1.15%  ' $false